$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3955.3635
$ws.Range("I62").Value = 2363.3
$ws.Range("K62").Value = 2363.3
$ws.Range("M62").Value = -1739.3

$ws.Range("H65").Value = 3955.3635
$ws.Range("I65").Value = 2363.3
$ws.Range("K65").Value = 11816.5
$ws.Range("M65").Value = -8696.5

$ws.Range("H70").Value = 1781.625
$ws.Range("I70").Value = 1542.8572
$ws.Range("J70").Value = 2237.4546
$ws.Range("K70").Value = 4628.571599999999
$ws.Range("L70").Value = 6712.3638
$ws.Range("M70").Value = -4358.571599999999
$ws.Range("N70").Value = -7252.3638

$ws.Range("H73").Value = 1781.625
$ws.Range("I73").Value = 1542.8572
$ws.Range("J73").Value = 2237.4546
$ws.Range("K73").Value = 4628.571599999999
$ws.Range("L73").Value = 6712.3638
$ws.Range("M73").Value = -3692.571599999999
$ws.Range("N73").Value = -8584.363799999999

$ws.Range("H111").Value = 2558.0833
$ws.Range("I111").Value = 2022.1111
$ws.Range("J111").Value = 4166
$ws.Range("K111").Value = 6066.3333
$ws.Range("L111").Value = 12498
$ws.Range("M111").Value = -2999.3333
$ws.Range("N111").Value = -18632

$ws.Range("H112").Value = 9705
$ws.Range("J112").Value = 2155.625
$ws.Range("L112").Value = 6466.875
$ws.Range("N112").Value = -8682.875

$ws.Range("H125").Value = 1233
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 1137.125
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 10234.125
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -15154.125

$ws.Range("H129").Value = 949.7646999999999
$ws.Range("I129").Value = 389
$ws.Range("J129").Value = 1095.1482
$ws.Range("K129").Value = 1167
$ws.Range("L129").Value = 3285.4446
$ws.Range("M129").Value = 3833
$ws.Range("N129").Value = -13285.4446

$ws.Range("H130").Value = 80500
$ws.Range("J130").Value = 80500
$ws.Range("L130").Value = 80500
$ws.Range("N130").Value = -90540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1782.3334
$ws.Range("I45").Value = 1788.7778
$ws.Range("K45").Value = 1788.7778
$ws.Range("M45").Value = -1411.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 39461.6
$ws.Range("J110").Value = 39461.6
$ws.Range("L110").Value = 39461.6
$ws.Range("N110").Value = -47641.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 591500.1
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 622105.4
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 622105.4
$ws.Range("M4").Value = -9888
$ws.Range("N4").Value = -622329.4

$ws.Range("H99").Value = 1337.4286
$ws.Range("I99").Value = 1337.4286
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1337.4286
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 160.5714
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 470.72
$ws.Range("I107").Value = 321.75
$ws.Range("J107").Value = 735.55554
$ws.Range("K107").Value = 321.75
$ws.Range("L107").Value = 735.55554
$ws.Range("M107").Value = 1598.25
$ws.Range("N107").Value = -4575.55554

$ws.Range("H122").Value = 6146.967
$ws.Range("I122").Value = 2573.4092
$ws.Range("K122").Value = 7720.2276
$ws.Range("M122").Value = -5270.2276

$ws.Range("H126").Value = 1337.4286
$ws.Range("I126").Value = 1337.4286
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4012.2858
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1542.2858
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 11060
$ws.Range("I120").Value = 13626.667
$ws.Range("J120").Value = 10326.667
$ws.Range("K120").Value = 40880.001
$ws.Range("L120").Value = 30980.001
$ws.Range("M120").Value = -36042.001
$ws.Range("N120").Value = -40656.001

$ws.Range("H137").Value = 63993.375
$ws.Range("I137").Value = 1157.8334
$ws.Range("K137").Value = 3473.5002
$ws.Range("M137").Value = 1626.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3359.6667
$ws.Range("I2").Value = 39.5
$ws.Range("K2").Value = 39.5
$ws.Range("M2").Value = 73.5

$ws.Range("H11").Value = 280000000
$ws.Range("I11").Value = 280000000
$ws.Range("K11").Value = 280000000
$ws.Range("M11").Value = -279999861

$ws.Range("H18").Value = 7505001
$ws.Range("J18").Value = 6666.6665
$ws.Range("L18").Value = 6666.6665
$ws.Range("N18").Value = -7252.6665

$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H126").Value = 2597.3333
$ws.Range("I126").Value = 1681.1428
$ws.Range("J126").Value = 3880
$ws.Range("K126").Value = 5043.428400000001
$ws.Range("L126").Value = 11640
$ws.Range("M126").Value = -2573.428400000001
$ws.Range("N126").Value = -16580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4999
$ws.Range("J2").Value = 4999
$ws.Range("L2").Value = 4999
$ws.Range("N2").Value = -5223

$ws.Range("H7").Value = 2231.2144
$ws.Range("I7").Value = 2231.2144
$ws.Range("K7").Value = 2231.2144
$ws.Range("M7").Value = -2119.2144

$ws.Range("H13").Value = 1500
$ws.Range("J13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("N13").Value = -1780

$ws.Range("H25").Value = 30251.75
$ws.Range("I25").Value = 11007
$ws.Range("K25").Value = 11007
$ws.Range("M25").Value = -10777

$ws.Range("H40").Value = 3090.353
$ws.Range("I40").Value = 3000.2
$ws.Range("K40").Value = 3000.2
$ws.Range("M40").Value = -2864.2

$ws.Range("H122").Value = 7112.5
$ws.Range("I122").Value = 6850
$ws.Range("K122").Value = 20550
$ws.Range("M122").Value = -18100

$ws.Range("H126").Value = 2231.2144
$ws.Range("I126").Value = 2231.2144
$ws.Range("K126").Value = 6693.6432
$ws.Range("M126").Value = -4223.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1621.2778
$ws.Range("I107").Value = 1074.9
$ws.Range("J107").Value = 2304.25
$ws.Range("K107").Value = 3224.7
$ws.Range("L107").Value = 6912.75
$ws.Range("M107").Value = -1304.7
$ws.Range("N107").Value = -10752.75

$ws.Range("H126").Value = 1375.5
$ws.Range("I126").Value = 1290
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3870
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1400
$ws.Range("N126").Value = -13940
